$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)
$shp = $s.Shapes.Item(36)
$tbl = $shp.Table

# --- Resize the first two table columns (EMU 1371983 -> 1410559 / 1333407,
#     i.e. points = EMU / 12700) ---
$tbl.Columns.Item(1).Width = 111.06763779527559
$tbl.Columns.Item(2).Width = 104.99267716535434

# --- Row 2 ("Linear Regression (all variables)"): RMSE cell 2040 -> 20 ---
$tbl.Cell(2,4).Shape.TextFrame.TextRange.Text = "20"

# --- Row 4 ("CART"): finish the cp value and fill in OSR2 / MAE / RMSE ---
$tbl.Cell(4,1).Shape.TextFrame.TextRange.Text = "CART (cp = 1.5e-05"
$tbl.Cell(4,2).Shape.TextFrame.TextRange.Text = "0.0410"
$tbl.Cell(4,3).Shape.TextFrame.TextRange.Text = "26.17"
$tbl.Cell(4,4).Shape.TextFrame.TextRange.Text = "2042"
